$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new "irrigation" column (column F)
# Copy the formatting of the preceding header cell (E1) so the new
# header picks up the same bold/font style used by the other headers.
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122) # xlPasteFormats
$ws.Cells.Item(1, 6).Value = "irrigation"

# Irrigation (water volume) values for rows 2-145, column F
$irrigation = @(
    60,90,120,60,90,120,60,90,120,60,90,120,60,90,120,60,90,120,
    50,75,100,50,75,100,50,75,100,50,75,100,50,75,100,50,75,100,
    40,60,80,40,60,80,40,60,80,40,60,80,40,60,80,40,60,80,
    40,60,80,40,60,80,40,60,80,40,60,80,40,60,80,40,60,80,
    30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,
    30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,
    30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,30,45,60,
    0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0
)

for ($i = 0; $i -lt $irrigation.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $irrigation[$i]
}

# Update the view: scroll down and move the active selection, as the
# author did while entering the new data near the bottom of the sheet.
$ws.Range("A115").Select()
$excel.ActiveWindow.ScrollRow = 115
$ws.Range("H133").Select()
